# Programme Done - Without JDOC
#
# 1. Rename product "Mælk" -> "Mælkebøtte" (row 2)
# 2. Nudge A2's formatting (theme color + default orientation) so it
#    carries its own distinct style, matching the author's re-save.
# 3. Add a new product row "Flute" (row 9) with Q1..Q4 figures, using the
#    Calibri font to set it apart like the other distinctly-styled rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: update the existing "Mælk" row ---------------------------------
$ws.Range("A2").Value = "Mælkebøtte"
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Orientation = 0

# --- 3: append the new "Flute" row -----------------------------------------
$ws.Cells.Item(9, 1).Value = "Flute"
$ws.Cells.Item(9, 2).Value = 800.0
$ws.Cells.Item(9, 3).Value = 200.0
$ws.Cells.Item(9, 4).Value = 100.0
$ws.Cells.Item(9, 5).Value = 10.0
$ws.Range("A9:E9").Font.Name = "Calibri"

Write-Host "Applied: renamed A2, added Flute row 9"
